$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Il15"
$ws.Cells.Item(2, 3).Value = "Il2rg"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 16.58023833333333
$ws.Cells.Item(2, 8).Value = 49.740715
$ws.Cells.Item(2, 9).Value = 0.63541025828417
$ws.Cells.Item(2, 10).Value = 0.63541025828417
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 22.84733166666667
$ws.Cells.Item(2, 14).Value = 68.541995
$ws.Cells.Item(2, 15).Value = 0.8188848782804395
$ws.Cells.Item(2, 16).Value = 0.8188848782804397
$ws.Cells.Item(2, 17).Value = 378.8142043140472
$ws.Cells.Item(2, 18).Value = 3409.327838826425
$ws.Cells.Item(2, 19).Value = 0.5203278520131752
$ws.Cells.Item(2, 20).Value = 0.5203278520131752

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Il15"
$ws.Cells.Item(3, 3).Value = "Il2rg"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 16.58023833333333
$ws.Cells.Item(3, 8).Value = 49.740715
$ws.Cells.Item(3, 9).Value = 0.63541025828417
$ws.Cells.Item(3, 10).Value = 0.63541025828417
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.041337
$ws.Cells.Item(3, 14).Value = 0.124011
$ws.Cells.Item(3, 15).Value = 0.00148158413889814
$ws.Cells.Item(3, 16).Value = 0.00148158413889814
$ws.Cells.Item(3, 17).Value = 0.685377311985
$ws.Cells.Item(3, 18).Value = 6.168395807865
$ws.Cells.Item(3, 19).Value = 0.0009414137603669965
$ws.Cells.Item(3, 20).Value = 0.0009414137603669965

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Il15"
$ws.Cells.Item(4, 3).Value = "Il2rg"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 16.58023833333333
$ws.Cells.Item(4, 8).Value = 49.740715
$ws.Cells.Item(4, 9).Value = 0.63541025828417
$ws.Cells.Item(4, 10).Value = 0.63541025828417
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 5.011873
$ws.Cells.Item(4, 14).Value = 15.035619
$ws.Cells.Item(4, 15).Value = 0.1796335375806623
$ws.Cells.Item(4, 16).Value = 0.1796335375806623
$ws.Cells.Item(4, 17).Value = 83.09804883639835
$ws.Cells.Item(4, 18).Value = 747.882439527585
$ws.Cells.Item(4, 19).Value = 0.1141409925106278
$ws.Cells.Item(4, 20).Value = 0.1141409925106278

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Il15"
$ws.Cells.Item(5, 3).Value = "Il2rg"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 8.574149
$ws.Cells.Item(5, 8).Value = 25.722447
$ws.Cells.Item(5, 9).Value = 0.3285901035393414
$ws.Cells.Item(5, 10).Value = 0.3285901035393414
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 22.84733166666667
$ws.Cells.Item(5, 14).Value = 68.541995
$ws.Cells.Item(5, 15).Value = 0.8188848782804395
$ws.Cells.Item(5, 16).Value = 0.8188848782804397
$ws.Cells.Item(5, 17).Value = 195.8964259624183
$ws.Cells.Item(5, 18).Value = 1763.067833661765
$ws.Cells.Item(5, 19).Value = 0.2690774669409706
$ws.Cells.Item(5, 20).Value = 0.2690774669409707

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Il15"
$ws.Cells.Item(6, 3).Value = "Il2rg"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 8.574149
$ws.Cells.Item(6, 8).Value = 25.722447
$ws.Cells.Item(6, 9).Value = 0.3285901035393414
$ws.Cells.Item(6, 10).Value = 0.3285901035393414
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.041337
$ws.Cells.Item(6, 14).Value = 0.124011
$ws.Cells.Item(6, 15).Value = 0.00148158413889814
$ws.Cells.Item(6, 16).Value = 0.00148158413889814
$ws.Cells.Item(6, 17).Value = 0.354429597213
$ws.Cells.Item(6, 18).Value = 3.189866374917
$ws.Cells.Item(6, 19).Value = 0.0004868338856027858
$ws.Cells.Item(6, 20).Value = 0.0004868338856027858

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Il15"
$ws.Cells.Item(7, 3).Value = "Il2rg"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 8.574149
$ws.Cells.Item(7, 8).Value = 25.722447
$ws.Cells.Item(7, 9).Value = 0.3285901035393414
$ws.Cells.Item(7, 10).Value = 0.3285901035393414
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 5.011873
$ws.Cells.Item(7, 14).Value = 15.035619
$ws.Cells.Item(7, 15).Value = 0.1796335375806623
$ws.Cells.Item(7, 16).Value = 0.1796335375806623
$ws.Cells.Item(7, 17).Value = 42.972545871077
$ws.Cells.Item(7, 18).Value = 386.752912839693
$ws.Cells.Item(7, 19).Value = 0.059025802712768
$ws.Cells.Item(7, 20).Value = 0.059025802712768

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Il15"
$ws.Cells.Item(8, 3).Value = "Il2rg"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.9393656666666668
$ws.Cells.Item(8, 8).Value = 2.818097
$ws.Cells.Item(8, 9).Value = 0.03599963817648871
$ws.Cells.Item(8, 10).Value = 0.03599963817648871
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 22.84733166666667
$ws.Cells.Item(8, 14).Value = 68.541995
$ws.Cells.Item(8, 15).Value = 0.8188848782804395
$ws.Cells.Item(8, 16).Value = 0.8188848782804397
$ws.Cells.Item(8, 17).Value = 21.46199894261278
$ws.Cells.Item(8, 18).Value = 193.157990483515
$ws.Cells.Item(8, 19).Value = 0.02947955932629382
$ws.Cells.Item(8, 20).Value = 0.02947955932629383

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Il15"
$ws.Cells.Item(9, 3).Value = "Il2rg"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.9393656666666668
$ws.Cells.Item(9, 8).Value = 2.818097
$ws.Cells.Item(9, 9).Value = 0.03599963817648871
$ws.Cells.Item(9, 10).Value = 0.03599963817648871
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.041337
$ws.Cells.Item(9, 14).Value = 0.124011
$ws.Cells.Item(9, 15).Value = 0.00148158413889814
$ws.Cells.Item(9, 16).Value = 0.00148158413889814
$ws.Cells.Item(9, 17).Value = 0.038830558563
$ws.Cells.Item(9, 18).Value = 0.349475027067
$ws.Cells.Item(9, 19).Value = 0.00005333649292835762
$ws.Cells.Item(9, 20).Value = 0.00005333649292835763

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Il15"
$ws.Cells.Item(10, 3).Value = "Il2rg"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.9393656666666668
$ws.Cells.Item(10, 8).Value = 2.818097
$ws.Cells.Item(10, 9).Value = 0.03599963817648871
$ws.Cells.Item(10, 10).Value = 0.03599963817648871
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 5.011873
$ws.Cells.Item(10, 14).Value = 15.035619
$ws.Cells.Item(10, 15).Value = 0.1796335375806623
$ws.Cells.Item(10, 16).Value = 0.1796335375806623
$ws.Cells.Item(10, 17).Value = 4.707981421893668
$ws.Cells.Item(10, 18).Value = 42.37183279704301
$ws.Cells.Item(10, 19).Value = 0.006466742357266528
$ws.Cells.Item(10, 20).Value = 0.00646674235726653
